# Add the new data row (row 7) reported after running the profit script on 2025-12-01.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A holds a date formatted as plain text (like the existing rows), so force
# text storage before assigning it — otherwise Excel auto-converts the
# "MM/DD/YYYY"-looking string into a real date serial number.
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "12/01/2025"
$ws.Cells.Item(7, 1).Style = "Normal"

$ws.Cells.Item(7, 2).Value = 12960.19
$ws.Cells.Item(7, 3).Value = 0.1667601265990996
$ws.Cells.Item(7, 4).Value = 0.8332398734009004
$ws.Cells.Item(7, 5).Value = -94.72
$ws.Cells.Item(7, 6).Value = -23.1
$ws.Cells.Item(7, 7).Value = -19667.81
$ws.Cells.Item(7, 8).Value = -64.55
$ws.Cells.Item(7, 9).Value = -641.61
$ws.Cells.Item(7, 10).Value = -22.89
